$wb = $excel.ActiveWorkbook

# 1. Update TS_Defs!Q6 from "p,t" to "t" (this makes the old "p,t" shared string unused)
$wsTSDefs = $wb.Worksheets.Item("TS_Defs")
$wsTSDefs.Range("Q6").Value = "t"

# 2. Add two new rows to the "process map" sheet
$wsProcMap = $wb.Worksheets.Item("process map")
$wsProcMap.Range("A24").Value = "old_new"
$wsProcMap.Range("B24").Value = "*"
$wsProcMap.Range("C24").Value = "new"
$wsProcMap.Range("A25").Value = "old_new"
$wsProcMap.Range("B25").Value = "ep*"
$wsProcMap.Range("C25").Value = "old"

# 3. Update selection on TS_Defs to A6
$wsTSDefs.Activate() | Out-Null
$wsTSDefs.Range("A6").Select() | Out-Null

# 4. Make "process map" the active sheet, and select C26
$wsProcMap.Activate() | Out-Null
$wsProcMap.Range("C26").Select() | Out-Null
